# Add a new worksheet "testClaimDashboard" after "testCancelclaim", modeled
# after the commit "Added dashboard claims testcase for claimns".
#
# Strategy: copy the last existing sheet (testCancelclaim) so the new sheet
# inherits the workbook's normal per-sheet defaults (row height, x14ac
# namespace, etc.), then rename it and replace its contents with the new
# claim-dashboard header/data row.

$wb = $excel.ActiveWorkbook

$srcSheet = $wb.Worksheets.Item("testCancelclaim")
$srcSheet.Copy([System.Reflection.Missing]::Value, $srcSheet)

$ws = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws.Name = "testClaimDashboard"

# Clear out the copied content before writing the new table.
$ws.Cells.ClearContents()

# Header row.
$ws.Range("A1").Value = "Username"
$ws.Range("B1").Value = "Password"
$ws.Range("C1").Value = "Header1"
$ws.Range("D1").Value = "Header2"
$ws.Range("E1").Value = "Header3"
$ws.Range("F1").Value = "Header4"

# Data row.
$ws.Range("A2").Value = "Admin"
$ws.Range("B2").Value = "admin123"
$ws.Range("C2").Value = "Submit Claim"
$ws.Range("D2").Value = "My Claims"
$ws.Range("E2").Value = "Employee Claims"
$ws.Range("F2").Value = "Assign Claim"

# Column widths for the widest header/value columns.
$ws.Range("C1:C2").ColumnWidth = 11.1796875
$ws.Range("G1:G2").ColumnWidth = 12.1796875

# Leave the selection the way the source workbook had it.
$ws.Range("F2").Select()

# The newly added sheet becomes the active tab (matches the source diff,
# where testCancelclaim's tabSelected moves to the new testClaimDashboard
# sheet and its own selection collapses to a single cell).
$srcSheet.Range("G15").Select()
$ws.Activate()
